$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.129.55"
$ws.Range("E2").Value = "  +7.82%  "

# Row 3
$ws.Range("D3").Value = "2.532.40"
$ws.Range("E3").Value = "  +7.97%  "

# Row 4
$ws.Range("E4").Value = "  +0.32%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.23%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +25.22%  "

# Row 8
$ws.Range("E8").Value = "  -0.78%  "

# Row 9
$ws.Range("D9").Value = "2.580.49"
$ws.Range("E9").Value = "  +10.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.28%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.95%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.82%  "

# Row 14
$ws.Range("D14").Value = "2.980.08"
$ws.Range("E14").Value = "  +8.06%  "

# Row 15
$ws.Range("D15").Value = "59.067.83"
$ws.Range("E15").Value = "  +7.50%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.94%  "

# Row 17
$ws.Range("E17").Value = "  +5.63%  "

# Row 18
$ws.Range("D18").Value = "2.571.04"
$ws.Range("E18").Value = "  +9.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "333.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.54%  "

# Row 21
$ws.Range("E21").Value = "  +8.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.91%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.414"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.70%  "

# Row 27
$ws.Range("D27").Value = "2.667.23"
$ws.Range("E27").Value = "  +8.78%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.34%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0826"
$ws.Range("E29").Value = "  +10.80%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.47%  "

# Row 34
$ws.Range("E34").Value = "  +7.59%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.88%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.18%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.51%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "290.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.57%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.101"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.622"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "  +7.59%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.990"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.89%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.755"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +21.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.86%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.98%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0235"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.43%  "

# Row 51
$ws.Range("D51").Value = "1.998.61"
$ws.Range("E51").Value = "  +12.68%  "
